$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 117.9639543333333
$ws.Range("H2").Value = 353.891863
$ws.Range("I2").Value = 0.2661690114309019
$ws.Range("J2").Value = 0.2661690114309019
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 117.044563
$ws.Range("N2").Value = 351.133689
$ws.Range("O2").Value = 0.3245365645427815
$ws.Range("P2").Value = 0.3245365645427815
$ws.Range("Q2").Value = 13807.03948469696
$ws.Range("R2").Value = 124263.3553622726
$ws.Range("S2").Value = 0.08638157655753324
$ws.Range("T2").Value = 0.08638157655753324

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 117.9639543333333
$ws.Range("H3").Value = 353.891863
$ws.Range("I3").Value = 0.2661690114309019
$ws.Range("J3").Value = 0.2661690114309019
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 101.5800373333333
$ws.Range("N3").Value = 304.740112
$ws.Range("O3").Value = 0.281657135515876
$ws.Range("P3").Value = 0.281657135515876
$ws.Range("Q3").Value = 11982.78288516763
$ws.Range("R3").Value = 107845.0459665087
$ws.Range("S3").Value = 0.0749684013227203
$ws.Range("T3").Value = 0.07496840132272029

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 117.9639543333333
$ws.Range("H4").Value = 353.891863
$ws.Range("I4").Value = 0.2661690114309019
$ws.Range("J4").Value = 0.2661690114309019
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 142.0267893333333
$ws.Range("N4").Value = 426.080368
$ws.Range("O4").Value = 0.3938062999413425
$ws.Range("P4").Value = 0.3938062999413425
$ws.Range("Q4").Value = 16754.04169102729
$ws.Range("R4").Value = 150786.3752192456
$ws.Range("S4").Value = 0.1048190335506484
$ws.Range("T4").Value = 0.1048190335506484

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 282.6413673333333
$ws.Range("H5").Value = 847.924102
$ws.Range("I5").Value = 0.6377403483780447
$ws.Range("J5").Value = 0.6377403483780446
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 117.044563
$ws.Range("N5").Value = 351.133689
$ws.Range("O5").Value = 0.3245365645427815
$ws.Range("P5").Value = 0.3245365645427815
$ws.Range("Q5").Value = 33081.63532525247
$ws.Range("R5").Value = 297734.7179272723
$ws.Range("S5").Value = 0.2069700617329272
$ws.Range("T5").Value = 0.2069700617329272

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 282.6413673333333
$ws.Range("H6").Value = 847.924102
$ws.Range("I6").Value = 0.6377403483780447
$ws.Range("J6").Value = 0.6377403483780446
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 101.5800373333333
$ws.Range("N6").Value = 304.740112
$ws.Range("O6").Value = 0.281657135515876
$ws.Range("P6").Value = 0.281657135515876
$ws.Range("Q6").Value = 28710.72064566438
$ws.Range("R6").Value = 258396.4858109794
$ws.Range("S6").Value = 0.1796241197270569
$ws.Range("T6").Value = 0.1796241197270569

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 282.6413673333333
$ws.Range("H7").Value = 847.924102
$ws.Range("I7").Value = 0.6377403483780447
$ws.Range("J7").Value = 0.6377403483780446
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 142.0267893333333
$ws.Range("N7").Value = 426.080368
$ws.Range("O7").Value = 0.3938062999413425
$ws.Range("P7").Value = 0.3938062999413425
$ws.Range("Q7").Value = 40142.64593513661
$ws.Range("R7").Value = 361283.8134162296
$ws.Range("S7").Value = 0.2511461669180606
$ws.Range("T7").Value = 0.2511461669180605

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 42.586595
$ws.Range("H8").Value = 127.759785
$ws.Range("I8").Value = 0.09609064019105341
$ws.Range("J8").Value = 0.09609064019105343
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 117.044563
$ws.Range("N8").Value = 351.133689
$ws.Range("O8").Value = 0.3245365645427815
$ws.Range("P8").Value = 0.3245365645427815
$ws.Range("Q8").Value = 4984.529401432985
$ws.Range("R8").Value = 44860.76461289686
$ws.Range("S8").Value = 0.031184926252321
$ws.Range("T8").Value = 0.031184926252321

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 42.586595
$ws.Range("H9").Value = 127.759785
$ws.Range("I9").Value = 0.09609064019105341
$ws.Range("J9").Value = 0.09609064019105343
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 101.5800373333333
$ws.Range("N9").Value = 304.740112
$ws.Range("O9").Value = 0.281657135515876
$ws.Range("P9").Value = 0.281657135515876
$ws.Range("Q9").Value = 4325.947909999546
$ws.Range("R9").Value = 38933.53118999592
$ws.Range("S9").Value = 0.02706461446609881
$ws.Range("T9").Value = 0.02706461446609881

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 42.586595
$ws.Range("H10").Value = 127.759785
$ws.Range("I10").Value = 0.09609064019105341
$ws.Range("J10").Value = 0.09609064019105343
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 142.0267893333333
$ws.Range("N10").Value = 426.080368
$ws.Range("O10").Value = 0.3938062999413425
$ws.Range("P10").Value = 0.3938062999413425
$ws.Range("Q10").Value = 6048.437356488987
$ws.Range("R10").Value = 54435.93620840088
$ws.Range("S10").Value = 0.0378410994726336
$ws.Range("T10").Value = 0.03784109947263361
